$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the three runs that make up
#      "***IMAGE imageGenerator(ara.logoUrl, " + "6, 2" + ")***"
#    (all three share identical run formatting) into the same text kept in
#    a single run.
# ---------------------------------------------------------------------------
$imgRng = $d.Content
$imgFound = $imgRng.Find.Execute("***IMAGE imageGenerator(ara.logoUrl, 6, 2)***", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "***IMAGE imageGenerator(ara.logoUrl, 6, 2)***", 2)

# ---------------------------------------------------------------------------
# 2) Invoice label fix: "Jura: " -> "Juros: " (refs #1538).
#    In the source edit this was produced by replacing just the "a" in the
#    middle of the word with "os", which leaves the surrounding "Jur" and
#    ": " text runs intact and creates a new run for the inserted "os".
#    We reproduce that exact run layout (3 runs, identical run formatting)
#    via InsertXML so the three runs stay distinct instead of being
#    recombined into a single run.
# ---------------------------------------------------------------------------
$juraRng = $d.Content
$juraFound = $juraRng.Find.Execute("Jura: ")

if ($juraFound) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:rPr></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Jur</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>os</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $juraRng.InsertXML($xml)
}
